# Update timetable worksheets (Section_A, Section_B) and the
# Elective_Coordination schedule to reflect the rebalanced slots.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Section_A
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Section_A")

$ws.Range("B2").Value = "Free"
$ws.Range("C2").Value = "Free"
$ws.Range("D2").Value = "CS251 (Elective)"
$ws.Range("E2").Value = "DA262"

$ws.Range("B3").Value = "CS307"
$ws.Range("C3").Value = "HS261 (Elective)"
$ws.Range("D3").Value = "CS307"
$ws.Range("E3").Value = "CS307"
$ws.Range("F3").Value = "Free"

$ws.Range("B5").Value = "DA262"
$ws.Range("C5").Value = "CS251 (Elective)"
$ws.Range("D5").Value = "DA261"
$ws.Range("F5").Value = "CS304"

$ws.Range("B6").Value = "HS261 (Tutorial)"
$ws.Range("F6").Value = "Free"

$ws.Range("B7").Value = "CS304"
$ws.Range("C7").Value = "CS304"
$ws.Range("D7").Value = "Free"
$ws.Range("F7").Value = "HS261 (Elective)"

$ws.Range("C8").Value = "Free"
$ws.Range("F8").Value = "CS251 (Tutorial)"

# ---------------------------------------------------------------
# Section_B
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Section_B")

$ws.Range("B2").Value = "Free"
$ws.Range("C2").Value = "CS304"
$ws.Range("D2").Value = "CS251 (Elective)"
$ws.Range("F2").Value = "Free"

$ws.Range("B3").Value = "Free"
$ws.Range("C3").Value = "HS261 (Elective)"
$ws.Range("D3").Value = "CS307"
$ws.Range("E3").Value = "CS304"
$ws.Range("F3").Value = "DA262"

$ws.Range("B5").Value = "CS304"
$ws.Range("C5").Value = "CS251 (Elective)"
$ws.Range("D5").Value = "Free"
$ws.Range("E5").Value = "CS307"
$ws.Range("F5").Value = "CS307"

$ws.Range("B6").Value = "HS261 (Tutorial)"
$ws.Range("C6").Value = "Free"
$ws.Range("E6").Value = "CS304 (Tutorial)"
$ws.Range("F6").Value = "Free"

$ws.Range("B7").Value = "DA262"
$ws.Range("C7").Value = "DA261"
$ws.Range("D7").Value = "Free"
$ws.Range("E7").Value = "Free"
$ws.Range("F7").Value = "HS261 (Elective)"

$ws.Range("C8").Value = "Free"
$ws.Range("F8").Value = "CS251 (Tutorial)"

# ---------------------------------------------------------------
# Elective_Coordination
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Elective_Coordination")

$ws.Range("C2").Value = "Fri"
$ws.Range("D2").Value = "15:30-17:00"

$ws.Range("C3").Value = "Tue"
$ws.Range("D3").Value = "10:30-12:00"

$ws.Range("C4").Value = "Mon"
$ws.Range("D4").Value = "14:30-15:30"

$ws.Range("C11").Value = "Wed"
$ws.Range("D11").Value = "09:00-10:30"

$ws.Range("C12").Value = "Tue"
$ws.Range("D12").Value = "13:00-14:30"

$ws.Range("D13").Value = "17:00-18:00"
